# Update "想去人数" (interest count) values in F column across sheets,
# matching the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 268
$ws.Range("F5").Value  = 328
$ws.Range("F6").Value  = 473
$ws.Range("F7").Value  = 2184
$ws.Range("F9").Value  = 58
$ws.Range("F10").Value = 1636
$ws.Range("F11").Value = 1636
$ws.Range("F12").Value = 1364
$ws.Range("F14").Value = 1412
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 580
$ws.Range("F18").Value = 162
$ws.Range("F20").Value = 7258
$ws.Range("F21").Value = 7994
$ws.Range("F22").Value = 48
$ws.Range("F24").Value = 201
$ws.Range("F27").Value = 91
$ws.Range("F28").Value = 218
$ws.Range("F29").Value = 264
$ws.Range("F35").Value = 1442
$ws.Range("F37").Value = 228
$ws.Range("F40").Value = 13
$ws.Range("F41").Value = 732
$ws.Range("F43").Value = 1364
$ws.Range("F44").Value = 344
$ws.Range("F45").Value = 246
$ws.Range("F46").Value = 197
$ws.Range("F49").Value = 161

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 299

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2623
$ws.Range("F4").Value = 283
$ws.Range("F5").Value = 138

# Sheet "全部类型" (All types - aggregate of the above)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 268
$ws.Range("F6").Value  = 138
$ws.Range("F7").Value  = 328
$ws.Range("F9").Value  = 473
$ws.Range("F10").Value = 2184
$ws.Range("F12").Value = 58
$ws.Range("F13").Value = 1636
$ws.Range("F14").Value = 1636
$ws.Range("F16").Value = 1412
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 580
$ws.Range("F20").Value = 162
$ws.Range("F24").Value = 7258
$ws.Range("F25").Value = 7994
$ws.Range("F26").Value = 48
$ws.Range("F28").Value = 91
$ws.Range("F31").Value = 1442
$ws.Range("F33").Value = 228
$ws.Range("F37").Value = 13
$ws.Range("F39").Value = 732
$ws.Range("F43").Value = 1364
$ws.Range("F44").Value = 344
$ws.Range("F45").Value = 246
$ws.Range("F46").Value = 197
$ws.Range("F49").Value = 299

$wb.Save()
